$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tweak a few of the date values first.
$ws.Range("C2").Value = "15"
$ws.Range("F2").Value = "18"
$ws.Range("C4").Value = "25"

# Insert a new column before column A; this shifts the old A:J data to B:K,
# carrying cell styles along with it.
$ws.Columns("A").EntireColumn.Insert()

# New first column: "City" header + constant city value for every data row.
$ws.Range("A1").Value = "City"
$ws.Range("A2").Value = "Cheese City Per Person, IL, US"
$ws.Range("A3").Value = "Cheese City Per Person, IL, US"
$ws.Range("A4").Value = "Cheese City Per Person, IL, US"
$ws.Range("A5").Value = "Cheese City Per Person, IL, US"

# New last column: "hotelname" header + a hotel name per booking row.
# Copy formatting from the neighbouring column (K) before writing values so
# the new cells share the same style records instead of minting new ones.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "hotelname"

$ws.Range("L2").Value = "Cheese Automation Testing Resort"
$ws.Range("L3").Value = "Chocolate Cheesecake Beach"

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = "Cheese Cake Island Hotel"

$ws.Range("L5").Value = "Cheese Automation Testing Resort"

$ws.Range("L5").Select()
